$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A=44004.3257520486; B="dralubarcellos@gmail.com"; C="Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III"; D=8; E=0; F=6; G=0; H=0; I="utipoa" },
    @{ A=44004.3372102662; B="fnagel@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III"; D=123; E=1; F=96; G=6; H=38; I="utipoa" },
    @{ A=44004.3381328241; B="fnagel@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III"; D=128; E=0; F=96; G=5; H=43; I="utipoa" },
    @{ A=44004.3435996528; B="cdalmora@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI PEDIATRICA - TIPO III"; D=13; E=0; F=12; G=3; H=0; I="utipoa" },
    @{ A=44004.3455832639; B="dralubarcellos@gmail.com"; C="Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III"; D=8; E=0; F=7; G=0; H=0; I="utipoa" },
    @{ A=44004.3843541782; B="fernanda.stringhi@maededeus.com.br"; C="Hospital Mãe de Deus - UTI ADULTO - TIPO I"; D=60; E=0; F=48; G=9; H=4; I="utipoa" },
    @{ A=44004.3897736227; B="francojw66@yahoo.com.br"; C="Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS"; D=4; E=0; F=4; G=0; H=0; I="utipoa" },
    @{ A=44004.3924785532; B="pedrocomerlato@gmail.com"; C="Hospital Independência - UTI ADULTO - TIPO II"; D=10; E=0; F=10; G=0; H=0; I="utipoa" },
    @{ A=44004.3996724421; B="fredikg@yahoo.com.br"; C="Hospital da Restinga - UTI ADULTO - TIPO II"; D=10; E=0; F=10; G=0; H=0; I="utipoa" },
    @{ A=44004.3997353472; B="fernanda.stringhi@maededeus.com.br"; C="Hospital Mãe de Deus - UTI ADULTO - TIPO I"; D=60; E=0; F=48; G=9; H=4; I="utipoa" },
    @{ A=44004.3997571644; B="leandra@portoalegre.rs.gov.br"; C="Hospital Materno Infantil Presidente Vargas - UTI PEDIATRICA - TIPO II"; D=12; E=0; F=6; G=0; H=0; I="utipoa" },
    @{ A=44004.4164537037; B="lufacchi@uol.com.be"; C="Hospital Vila Nova - UTI ADULTO - TIPO II"; D=20; E=0; F=20; G=2; H=0; I="utipoa" },
    @{ A=44004.4196158218; B="analise.medina@divinaprovidencia.org.br"; C="Hospital Divina Providência - UTI ADULTO - TIPO II"; D=20; E=0; F=18; G=2; H=2; I="utipoa" },
    @{ A=44004.4255394676; B="taianivargas@hotmail.com"; C="Hospital Nossa Senhora da Conceição - UTI ADULTO - TIPO III"; D=75; E=0; F=70; G=1; H=27; I="utipoa" },
    @{ A=44004.4264949884; B="braun.luiz@gmail.com"; C="Hospital Nossa Senhora da Conceição - UTI PEDIATRICA - TIPO II"; D=18; E=4; F=6; G=2; H=0; I="utipoa" },
    @{ A=44004.4270557292; B="taianivargas@hotmail.com"; C="Hospital Nossa Senhora da Conceição - UTI ADULTO - TIPO III"; D=75; E=0; F=69; G=1; H=25; I="utipoa" },
    @{ A=44004.435489456; B="roseuti@gmail.com"; C="Hospital Moinhos de Vento - UTI ADULTO - TIPO III"; D=56; E=0; F=39; G=0; H=11; I="utipoa" },
    @{ A=44004.4513226968; B="smarcos@ghc.com.br"; C="Hospital Femina - UTI ADULTO - TIPO II"; D=6; E=0; F=3; G=0; H=0; I="utipoa" },
    @{ A=44004.4762848495; B="lauren.ghion@santacasa.org.br"; C="Complexo Hospitalar Santa Casa - UTI ADULTO - TIPO III"; D=87; E=0; F=70; G=3; H=6; I="utipoa" },
    @{ A=44004.4773660069; B="lauren.ghion@santacasa.org.br"; C="Complexo Hospitalar Santa Casa - UTI PEDIATRICA - TIPO III"; D=37; E=3; F=33; G=0; H=0; I="utipoa" },
    @{ A=44004.5038592245; B="fernanda.stringhi@maededeus.com.br"; C="Hospital Mãe de Deus - UTI ADULTO - TIPO I"; D=60; E=0; F=48; G=9; H=4; I="utipoa" },
    @{ A=44004.5080246065; B="fernanda.stringhi@maededeus.com.br"; C="Hospital Mãe de Deus - UTI ADULTO - TIPO I"; D=60; E=0; F=48; G=9; H=4; I="utipoa" },
    @{ A=44004.5566112037; B="ccih@hpa.org.br"; C="Hospital Porto Alegre - UTI ADULTO - TIPO II"; D=7; E=0; F=6; G=0; H=1; I="utipoa" },
    @{ A=44004.5597599537; B="vivianmed83@yahoo.com.br"; C="Hospital de Pronto Socorro de Porto Alegre - UTI ADULTO - TIPO II"; D=20; E=0; F=13; G=0; H=0; I="utipoa" },
    @{ A=44004.6316032176; B="joao.krauzer@hmv.org.br"; C="Hospital Moinhos de Vento - UTI PEDIATRICA - TIPO III"; D=11; E=0; F=5; G=0; H=0; I="utipoa" },
    @{ A=44004.6769332176; B="renatocvaz@hotmail.com"; C="Instituto de Cardiologia - UTI ADULTO - TIPO III"; D=51; E=3; F=30; G=1; H=2; I="utipoa" },
    @{ A=44004.6816339699; B="marcosboniatti@gmail.com"; C="Hospital Cristo Redentor - UTI ADULTO - TIPO III"; D=39; E=0; F=32; G=0; H=0; I="utipoa" },
    @{ A=44004.6919118866; B="francojw66@yahoo.com.br"; C="Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS"; D=4; E=0; F=4; G=0; H=0; I="utipoa" },
    @{ A=44004.7087064583; B="lufacchi@uol.com.br"; C="Hospital Vila Nova - UTI ADULTO - TIPO II"; D=20; E=0; F=20; G=1; H=1; I="utipoa" },
    @{ A=44004.7626814583; B="renatafarinon@yahoo.com.br"; C="Hospital São Lucas - UTI ADULTO - TIPO III"; D=59; E=0; F=41; G=1; H=4; I="utipoa" },
    @{ A=44004.7628308565; B="renatafarinon@yahoo.com.br"; C="Hospital Santa Ana - UTI ADULTO - TIPO II"; D=10; E=0; F=3; G=0; H=0; I="utipoa" },
    @{ A=44004.7827978935; B="smarcos@ghc.com.br"; C="Hospital Femina - UTI ADULTO - TIPO II"; D=6; E=0; F=3; G=0; H=0; I="utipoa" },
    @{ A=44005.2764640394; B="fnagel@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III"; D=128; E=0; F=96; G=5; H=43; I="utipoa" },
    @{ A=44005.277845463; B="fnagel@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI ADULTO - TIPO III"; D=128; E=0; F=96; G=5; H=45; I="utipoa" },
    @{ A=44005.3223271181; B="andre.machado@hed.com.br"; C="Hospital Ernesto Dorenelles - UTI ADULTO - TIPO III"; D=40; E=0; F=29; G=5; H=3; I="utipoa" },
    @{ A=44005.3237513889; B="analuizafilipini@gmail.com"; C="Hospital São Lucas - UTI ADULTO - TIPO III"; D=59; E=0; F=47; G=6; H=4; I="utipoa" },
    @{ A=44005.3253003241; B="dralubarcellos@gmail.com"; C="Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III"; D=8; E=0; F=4; G=0; H=0; I="utipoa" },
    @{ A=44005.3351635648; B="cdalmora@hcpa.edu.br"; C="Hospital de Clínicas de Porto Alegre - UTI PEDIATRICA - TIPO III"; D=13; E=0; F=10; G=1; H=0; I="utipoa" },
    @{ A=44005.379507419; B="lufacchi@uol.com.br"; C="Hospital Vila Nova - UTI ADULTO - TIPO II"; D=20; E=0; F=19; G=1; H=0; I="utipoa" },
    @{ A=44005.3822382176; B="leandra@portoalegre.rs.gov.br"; C="Hospital Materno Infantil Presidente Vargas - UTI PEDIATRICA - TIPO II"; D=12; E=0; F=3; G=0; H=0; I="utipoa" },
    @{ A=44005.3836240856; B="analise.medina@divinaprovidencia.org.br"; C="Hospital Divina Providência - UTI ADULTO - TIPO II"; D=20; E=0; F=18; G=1; H=2; I="utipoa" }
)

$startRow = 2744
$lastExistingRow = $startRow - 1

$r = $startRow
foreach ($row in $newRows) {
    # Copy the date-time style/number format from the template cell in column A
    $ws.Range("A$lastExistingRow").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I

    $r = $r + 1
}

Write-Output "Last row now:"
Write-Output $ws.UsedRange.Rows.Count